$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 63
$ws.Range("H2").Value = 63

$ws.Range("E3").Value = 37

$ws.Range("F5").Value = 89
$ws.Range("H5").Value = 89

$ws.Range("F6").Value = 31
$ws.Range("H6").Value = 31

$ws.Range("E10").Value = 482
$ws.Range("F10").Value = 246
$ws.Range("H10").Value = 246

$ws.Range("E11").Value = 324
$ws.Range("F11").Value = 184
$ws.Range("H11").Value = 184

$ws.Range("E12").Value = 472
$ws.Range("F12").Value = 259
$ws.Range("H12").Value = 259

$ws.Range("E13").Value = 123
$ws.Range("F13").Value = 66
$ws.Range("H13").Value = 66

$ws.Range("F14").Value = 63
$ws.Range("H14").Value = 63

$ws.Range("F15").Value = 67
$ws.Range("H15").Value = 67

$ws.Range("F16").Value = 96
$ws.Range("H16").Value = 96

$ws.Range("E17").Value = 91
$ws.Range("F17").Value = 48
$ws.Range("H17").Value = 48

$ws.Range("E20").Value = 84
$ws.Range("F20").Value = 29
$ws.Range("H20").Value = 29

$ws.Range("E22").Value = 161
$ws.Range("F22").Value = 86
$ws.Range("H22").Value = 86

$ws.Range("F23").Value = 88
$ws.Range("H23").Value = 88

$ws.Range("F24").Value = 110
$ws.Range("H24").Value = 110

$ws.Range("E25").Value = 251
$ws.Range("F25").Value = 124
$ws.Range("H25").Value = 124

$ws.Range("F26").Value = 92
$ws.Range("H26").Value = 92

$ws.Range("F27").Value = 154
$ws.Range("H27").Value = 154

$ws.Range("F28").Value = 72
$ws.Range("H28").Value = 72

$ws.Range("E29").Value = 163
$ws.Range("F29").Value = 92
$ws.Range("H29").Value = 92

$ws.Range("F30").Value = 118
$ws.Range("H30").Value = 118

$ws.Range("F32").Value = 105
$ws.Range("H32").Value = 105

$ws.Range("E33").Value = 276

$ws.Range("F34").Value = 135
$ws.Range("H34").Value = 135

$ws.Range("F35").Value = 89
$ws.Range("H35").Value = 89

$ws.Range("E37").Value = 150

$ws.Range("E38").Value = 88

$ws.Range("F39").Value = 86
$ws.Range("H39").Value = 86

$ws.Range("E40").Value = 250
$ws.Range("F40").Value = 119
$ws.Range("H40").Value = 119

$ws.Range("F42").Value = 194
$ws.Range("H42").Value = 194

$ws.Range("E43").Value = 110

$ws.Range("E44").Value = 300
$ws.Range("F44").Value = 152
$ws.Range("H44").Value = 152

$ws.Range("F46").Value = 168
$ws.Range("H46").Value = 168

$ws.Range("E47").Value = 429
$ws.Range("F47").Value = 212
$ws.Range("H47").Value = 212

$ws.Range("F48").Value = 85
$ws.Range("H48").Value = 85

$ws.Range("F49").Value = 119
$ws.Range("H49").Value = 119

$ws.Range("E50").Value = 235
$ws.Range("F50").Value = 109
$ws.Range("H50").Value = 109

$ws.Range("F51").Value = 95
$ws.Range("H51").Value = 95

